$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()

# Copy number formatting from column E (post-shift) into new column D
$srcRange = $ws.Range("D7:D102").Offset(0,1)
$srcRange.Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "done"
